$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (currently sitting at the
#    end of the "en github crear llave ssh y gpg keys" paragraph).
#    We'll re-create it later at its new location.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2. The paragraph that used to read
#      "git branch ver cuantas ramas tenemos (...)"
#    becomes the new first line:
#      "git push u origin master subir los cambios al repositorio"
# ------------------------------------------------------------------
$oldBranchText = "git branch ver cuantas ramas tenemos (cuando le damos a esto git no crea una rama temporal para viajar en el tiempo despues de los cambios que hicimos)"
$newPushText   = "git push u origin master subir los cambios al repositorio"

$null = $d.Content.Find.Execute($oldBranchText, $true, $false, $false, $false, $false, $true, 1, $false, $newPushText, 2)

# ------------------------------------------------------------------
# 3. Re-insert the original "git branch ver..." paragraph right after
#    the new "git push..." paragraph (so the old content re-appears,
#    just shifted one paragraph down).
# ------------------------------------------------------------------
$pushPara = $d.Paragraphs.Item(17)
$pushPara.Range.InsertParagraphAfter()
$branchPara = $d.Paragraphs.Item(18)
$branchPara.Range.Text = $oldBranchText

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark right after the new "git push..."
#    text (collapsed bookmark, immediately before the paragraph mark).
#    Directly creating a collapsed bookmark exactly at a paragraph's
#    content-end confuses Bookmarks.Add, so we insert a throw-away
#    character past the insertion point first, anchor the bookmark
#    before it, then remove the throw-away character again.
# ------------------------------------------------------------------
$pushPara = $d.Paragraphs.Item(17)
$contentEnd = $pushPara.Range.End - 1
$tmp = $d.Range($contentEnd, $contentEnd)
$tmp.InsertAfter("X")
$bmRange = $d.Range($contentEnd, $contentEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$killRange = $d.Range($contentEnd, $contentEnd + 1)
$killRange.Delete()

# ------------------------------------------------------------------
# 5. Extend the "en github crear llave ssh y gpg keys" paragraph text.
# ------------------------------------------------------------------
$oldKeysText = "en github crear llave ssh y gpg keys"
$newKeysText = "en github crear llave ssh y gpg keys para hacer la conexion mas segura al mmoento de subir el codigo"
$null = $d.Content.Find.Execute($oldKeysText, $true, $false, $false, $false, $false, $true, 1, $false, $newKeysText, 2)

# ------------------------------------------------------------------
# 6. Add the four new paragraphs right after it.
# ------------------------------------------------------------------
$keysPara = $d.Paragraphs.Item(24)

$keysPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item(25)
$p1.Range.Text = "y se puede tener varios key en cualquier pc trabjko u oficina "

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(26)
$p2.Range.Text = "crear proyecto "

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(27)
$p3.Range.Text = "en github"

# Final paragraph is a genuinely empty one (no run at all in the source
# XML), so rather than InsertParagraphAfter (which leaves a stray empty
# run behind) we split the paragraph mark directly.
$tail = $p3.Range.Duplicate
$tail.Collapse(0)
$tail.InsertAfter([char]13)

Write-Output "edit complete"
